$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 01:39"

# Swap rows 213 (Montserrat) and 214 (Islas Malvinas): Islas Malvinas moves up
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Updated country statistics
$ws.Range("B4").Value = 5475408
$ws.Range("C4").Value = 59742
$ws.Range("D4").Value = 2868271
$ws.Range("E4").Value = 2435660
$ws.Range("G4").Value = 1062
$ws.Range("H4").Value = 171477
$ws.Range("B5").Value = 3278895
$ws.Range("C5").Value = 49274
$ws.Range("E5").Value = 788022
$ws.Range("G5").Value = 1007
$ws.Range("H5").Value = 106571
$ws.Range("B11").Value = 445111
$ws.Range("C11").Value = 11306
$ws.Range("D11").Value = 261296
$ws.Range("E11").Value = 169323
$ws.Range("G11").Value = 347
$ws.Range("H11").Value = 14492
$ws.Range("B18").Value = 282437
$ws.Range("C18").Value = 6365
$ws.Range("E18").Value = 77905
$ws.Range("G18").Value = 165
$ws.Range("H18").Value = 5527
$ws.Range("B39").Value = 79402
$ws.Range("C39").Value = 956
$ws.Range("D39").Value = 52886
$ws.Range("E39").Value = 24782
$ws.Range("G39").Value = 12
$ws.Range("H39").Value = 1734
$ws.Range("B50").Value = 52217
$ws.Range("C50").Value = 1070
$ws.Range("D50").Value = 37479
$ws.Range("E50").Value = 13665
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 1073
$ws.Range("B52").Value = 48445
$ws.Range("C52").Value = 329
$ws.Range("D52").Value = 35998
$ws.Range("E52").Value = 11474
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 973
$ws.Range("B63").Value = 31381
$ws.Range("C63").Value = 1012
$ws.Range("D63").Value = 21580
$ws.Range("E63").Value = 9535
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 266
$ws.Range("B74").Value = 19693
$ws.Range("C74").Value = 292
$ws.Range("E74").Value = 5568
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 394
$ws.Range("B86").Value = 9908
$ws.Range("C86").Value = 57
$ws.Range("E86").Value = 790
$ws.Range("B98").Value = 7405
$ws.Range("C98").Value = 37
$ws.Range("D98").Value = 6500
$ws.Range("E98").Value = 783
$ws.Range("B113").Value = 3930
$ws.Range("C113").Value = 73
$ws.Range("D113").Value = 2752
$ws.Range("E113").Value = 1105
$ws.Range("B119").Value = 3229
$ws.Range("C119").Value = 55
$ws.Range("D119").Value = 2547
$ws.Range("E119").Value = 593
$ws.Range("B142").Value = 1421
$ws.Range("C142").Value = 12
$ws.Range("D142").Value = 1182
$ws.Range("E142").Value = 201
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 38
$ws.Range("B154").Value = 1119
$ws.Range("C154").Value = 30
$ws.Range("E154").Value = 964
$ws.Range("G154").Value = 2
$ws.Range("H154").Value = 17
$ws.Range("B159").Value = 930
$ws.Range("C159").Value = 19
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 22
$ws.Range("B168").Value = 426
$ws.Range("C168").Value = 22
$ws.Range("E168").Value = 277
$ws.Range("G168").Value = 2
$ws.Range("H168").Value = 10
